$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 699.3333
$ws.Range("I31").Value = 700
$ws.Range("J31").Value = 698
$ws.Range("K31").Value = 2100
$ws.Range("L31").Value = 2094
$ws.Range("M31").Value = -1870
$ws.Range("N31").Value = -2554
$ws.Range("H38").Value = 427.7143
$ws.Range("I38").Value = 110.333336
$ws.Range("K38").Value = 331.000008
$ws.Range("M38").Value = 40.99999200000002
$ws.Range("H39").Value = 228.28572
$ws.Range("I39").Value = 66.22221999999999
$ws.Range("J39").Value = 520
$ws.Range("K39").Value = 198.66666
$ws.Range("L39").Value = 1560
$ws.Range("M39").Value = 97.33334000000002
$ws.Range("N39").Value = -2152
$ws.Range("H40").Value = 3073.5454
$ws.Range("I40").Value = 15000
$ws.Range("J40").Value = 1880.9
$ws.Range("K40").Value = 15000
$ws.Range("L40").Value = 1880.9
$ws.Range("M40").Value = -14825
$ws.Range("N40").Value = -2230.9
$ws.Range("H112").Value = 4459.091
$ws.Range("I112").Value = 14369.8
$ws.Range("J112").Value = 1544.1765
$ws.Range("K112").Value = 43109.39999999999
$ws.Range("L112").Value = 4632.529500000001
$ws.Range("M112").Value = -42001.39999999999
$ws.Range("N112").Value = -6848.529500000001
$ws.Range("H129").Value = 971.8039
$ws.Range("I129").Value = 321.75
$ws.Range("J129").Value = 1027.1277
$ws.Range("K129").Value = 965.25
$ws.Range("L129").Value = 3081.3831
$ws.Range("M129").Value = 4034.75
$ws.Range("N129").Value = -13081.3831
$ws.Range("H137").Value = 3012.5938
$ws.Range("I137").Value = 3015.5
$ws.Range("K137").Value = 9046.5
$ws.Range("M137").Value = -6496.5
$ws.Range("H138").Value = 11580.703
$ws.Range("I138").Value = 3920.8572
$ws.Range("J138").Value = 13368
$ws.Range("K138").Value = 11762.5716
$ws.Range("L138").Value = 40104
$ws.Range("M138").Value = -6622.571599999999
$ws.Range("N138").Value = -50384

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8140.879
$ws.Range("I61").Value = 6161.8
$ws.Range("J61").Value = 11185.615
$ws.Range("K61").Value = 6161.8
$ws.Range("L61").Value = 11185.615
$ws.Range("M61").Value = -5949.8
$ws.Range("N61").Value = -11609.615
$ws.Range("H74").Value = 2212.0908
$ws.Range("I74").Value = 1987.8846
$ws.Range("K74").Value = 1987.8846
$ws.Range("M74").Value = -1113.8846
$ws.Range("H77").Value = 2212.0908
$ws.Range("I77").Value = 1987.8846
$ws.Range("K77").Value = 9939.423000000001
$ws.Range("M77").Value = -5571.423000000001
$ws.Range("H122").Value = 2136.1875
$ws.Range("I122").Value = 2168
$ws.Range("K122").Value = 6504
$ws.Range("M122").Value = -4054
$ws.Range("H136").Value = 8140.879
$ws.Range("I136").Value = 6161.8
$ws.Range("J136").Value = 11185.615
$ws.Range("K136").Value = 18485.4
$ws.Range("L136").Value = 33556.845
$ws.Range("M136").Value = -15935.4
$ws.Range("N136").Value = -38656.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3546
$ws.Range("I134").Value = 3266.6155
$ws.Range("K134").Value = 9799.8465
$ws.Range("M134").Value = -7264.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9996.275
$ws.Range("I31").Value = 10882.477
$ws.Range("J31").Value = 7670
$ws.Range("K31").Value = 10882.477
$ws.Range("L31").Value = 7670
$ws.Range("M31").Value = -10587.477
$ws.Range("N31").Value = -8260
$ws.Range("H34").Value = 9996.275
$ws.Range("I34").Value = 10882.477
$ws.Range("J34").Value = 7670
$ws.Range("K34").Value = 10882.477
$ws.Range("L34").Value = 7670
$ws.Range("M34").Value = -10680.477
$ws.Range("N34").Value = -8074
$ws.Range("H58").Value = 3498600
$ws.Range("I58").Value = 4786076
$ws.Range("K58").Value = 4786076
$ws.Range("M58").Value = -4785873
$ws.Range("H132").Value = 5372.2163
$ws.Range("I132").Value = 5620
$ws.Range("J132").Value = 4474
$ws.Range("K132").Value = 16860
$ws.Range("L132").Value = 13422
$ws.Range("M132").Value = -14330
$ws.Range("N132").Value = -18482
$ws.Range("H134").Value = 3082.8965
$ws.Range("I134").Value = 2477.4736
$ws.Range("J134").Value = 4233.2
$ws.Range("K134").Value = 7432.4208
$ws.Range("L134").Value = 12699.6
$ws.Range("M134").Value = -4897.4208
$ws.Range("N134").Value = -17769.6
$ws.Range("H136").Value = 3498600
$ws.Range("I136").Value = 4786076
$ws.Range("K136").Value = 14358228
$ws.Range("M136").Value = -14355678

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 38461770
$ws.Range("I12").Value = 166667020
$ws.Range("J12").Value = 193.15
$ws.Range("K12").Value = 500001060
$ws.Range("L12").Value = 579.45
$ws.Range("M12").Value = -500000887
$ws.Range("N12").Value = -925.45
$ws.Range("H48").Value = 3498
$ws.Range("J48").Value = 3498
$ws.Range("L48").Value = 10494
$ws.Range("N48").Value = -10994
$ws.Range("H60").Value = 370.75
$ws.Range("I60").Value = 212.66667
$ws.Range("J60").Value = 845
$ws.Range("K60").Value = 638.00001
$ws.Range("L60").Value = 2535
$ws.Range("M60").Value = -387.00001
$ws.Range("N60").Value = -3037
$ws.Range("H104").Value = 2805.9167
$ws.Range("J104").Value = 3279.1
$ws.Range("L104").Value = 9837.299999999999
$ws.Range("N104").Value = -15079.3
$ws.Range("H121").Value = 782.3333
$ws.Range("I121").Value = 576.6667
$ws.Range("J121").Value = 811.7143
$ws.Range("K121").Value = 1730.0001
$ws.Range("L121").Value = 2435.1429
$ws.Range("M121").Value = -420.0001
$ws.Range("N121").Value = -5055.1429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2506.205
$ws.Range("I102").Value = 1846.7273
$ws.Range("K102").Value = 1846.7273
$ws.Range("M102").Value = -224.7273
$ws.Range("H126").Value = 2611.739
$ws.Range("I126").Value = 1814.1666
$ws.Range("K126").Value = 5442.4998
$ws.Range("M126").Value = -2972.4998

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2746.8518
$ws.Range("I7").Value = 2974.7646
$ws.Range("K7").Value = 2974.7646
$ws.Range("M7").Value = -2862.7646
$ws.Range("H40").Value = 3808.9285
$ws.Range("I40").Value = 2883.3333
$ws.Range("J40").Value = 4503.125
$ws.Range("K40").Value = 2883.3333
$ws.Range("L40").Value = 4503.125
$ws.Range("M40").Value = -2747.3333
$ws.Range("N40").Value = -4775.125
$ws.Range("H122").Value = 6968.548
$ws.Range("I122").Value = 6513.357
$ws.Range("K122").Value = 19540.071
$ws.Range("M122").Value = -17090.071
$ws.Range("H126").Value = 2746.8518
$ws.Range("I126").Value = 2974.7646
$ws.Range("K126").Value = 8924.293799999999
$ws.Range("M126").Value = -6454.293799999999
$ws.Range("H132").Value = 5895.7144
$ws.Range("I132").Value = 5734.5557
$ws.Range("J132").Value = 6185.8
$ws.Range("K132").Value = 17203.6671
$ws.Range("L132").Value = 18557.4
$ws.Range("M132").Value = -14673.6671
$ws.Range("N132").Value = -23617.4
$ws.Range("H136").Value = 5414.421
$ws.Range("I136").Value = 3432.2354
$ws.Range("J136").Value = 7019.048
$ws.Range("K136").Value = 10296.7062
$ws.Range("L136").Value = 21057.144
$ws.Range("M136").Value = -7746.706200000001
$ws.Range("N136").Value = -26157.144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 741
$ws.Range("I113").Value = 403
$ws.Range("J113").Value = 1022.6667
$ws.Range("K113").Value = 1209
$ws.Range("L113").Value = 3068.0001
$ws.Range("M113").Value = 961
$ws.Range("N113").Value = -7408.0001
$ws.Range("H132").Value = 4151.5
$ws.Range("I132").Value = 3982
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 11946
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -9416
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 3855.946
$ws.Range("I136").Value = 2897.7896
$ws.Range("J136").Value = 4867.3335
$ws.Range("K136").Value = 8693.3688
$ws.Range("L136").Value = 14602.0005
$ws.Range("M136").Value = -6143.3688
$ws.Range("N136").Value = -19702.0005
